# Apply the maze grid update + selection change for spiral_hole2.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New maze values (0/1 grid), rows 1-10, columns A-J
$values = @(
  @(1,0,1,1,1,1,1,1,1,1),
  @(1,0,1,0,0,0,0,0,0,1),
  @(1,0,1,0,1,1,1,1,0,1),
  @(1,0,1,0,1,0,0,1,0,1),
  @(1,0,1,0,1,1,1,1,0,1),
  @(1,0,1,0,0,1,0,1,0,1),
  @(1,0,1,1,0,1,0,1,0,1),
  @(1,0,1,1,1,1,0,1,0,1),
  @(1,0,0,0,0,0,0,1,0,1),
  @(1,1,1,1,1,1,1,1,0,1)
)

for ($r = 0; $r -lt 10; $r++) {
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

# Update the active cell selection to match the author's edit
$ws.Range("V7").Select()
